$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2000
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").Value = ""

$ws.Range("H125").Value = 898.3333
$ws.Range("J125").Value = 530
$ws.Range("L125").Value = 4770
$ws.Range("N125").Value = -9690

$ws.Range("H138").Value = 9997.75
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 9997.75
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 29993.25
$ws.Range("M138").Value = ""
$ws.Range("N138").Value = -40273.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 5088
$ws.Range("I36").Value = 6684
$ws.Range("J36").Value = 300
$ws.Range("K36").Value = 6684
$ws.Range("L36").Value = 300
$ws.Range("M36").Value = -6338
$ws.Range("N36").Value = -992

$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").Value = ""

$ws.Range("H50").Value = 13653.714
$ws.Range("I50").Value = 5147.75
$ws.Range("J50").Value = 24995
$ws.Range("K50").Value = 5147.75
$ws.Range("L50").Value = 24995
$ws.Range("M50").Value = -4433.75
$ws.Range("N50").Value = -26423

$ws.Range("H102").Value = 6166.6665
$ws.Range("I102").Value = 5750
$ws.Range("K102").Value = 5750
$ws.Range("M102").Value = -4128

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").Value = ""

$ws.Range("H122").Value = 4235.875
$ws.Range("I122").Value = 4235.875
$ws.Range("K122").Value = 12707.625
$ws.Range("M122").Value = -10257.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 5333
$ws.Range("I107").Value = 4499.5
$ws.Range("K107").Value = 4499.5
$ws.Range("M107").Value = -2579.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 3288.7144
$ws.Range("J35").Value = 6516.3335
$ws.Range("L35").Value = 6516.3335
$ws.Range("N35").Value = -7104.3335

$ws.Range("H36").Value = 3969.7144
$ws.Range("J36").Value = 2153
$ws.Range("L36").Value = 2153
$ws.Range("N36").Value = -2929

$ws.Range("H40").Value = 3969.7144
$ws.Range("J40").Value = 2153
$ws.Range("L40").Value = 2153
$ws.Range("N40").Value = -2473

$ws.Range("H44").Value = 29831.166
$ws.Range("I44").Value = 29746.75
$ws.Range("J44").Value = 30000
$ws.Range("K44").Value = 29746.75
$ws.Range("L44").Value = 30000
$ws.Range("M44").Value = -29304.75
$ws.Range("N44").Value = -30884

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = ""
$ws.Range("N86").Value = ""

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = ""
$ws.Range("N89").Value = ""

$ws.Range("H132").Value = 2137.2727
$ws.Range("I132").Value = 2001
$ws.Range("K132").Value = 6003
$ws.Range("M132").Value = -3473

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H111").Value = 150
$ws.Range("I111").Value = 150
$ws.Range("K111").Value = 450
$ws.Range("M111").Value = 2617

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 1666.6666
$ws.Range("J43").Value = 2000
$ws.Range("L43").Value = 2000
$ws.Range("N43").Value = -2302

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2019.3636
$ws.Range("I22").Value = 960
$ws.Range("J22").Value = 2902.1667
$ws.Range("K22").Value = 960
$ws.Range("L22").Value = 2902.1667
$ws.Range("M22").Value = -665
$ws.Range("N22").Value = -3492.1667

$ws.Range("H26").Value = 4383.3335
$ws.Range("I26").Value = 3250
$ws.Range("J26").Value = 4950
$ws.Range("K26").Value = 3250
$ws.Range("L26").Value = 4950
$ws.Range("M26").Value = -2955
$ws.Range("N26").Value = -5540

$ws.Range("H27").Value = 2019.3636
$ws.Range("I27").Value = 960
$ws.Range("J27").Value = 2902.1667
$ws.Range("K27").Value = 960
$ws.Range("L27").Value = 2902.1667
$ws.Range("M27").Value = -853
$ws.Range("N27").Value = -3116.1667

$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").Value = ""

$ws.Range("H31").Value = 3241
$ws.Range("J31").Value = 3729.1667
$ws.Range("L31").Value = 3729.1667
$ws.Range("N31").Value = -4225.1667

$ws.Range("H46").Value = 287810.44
$ws.Range("I46").Value = 1000750
$ws.Range("J46").Value = 2634.6
$ws.Range("K46").Value = 1000750
$ws.Range("L46").Value = 2634.6
$ws.Range("M46").Value = -1000562
$ws.Range("N46").Value = -3010.6

$ws.Range("H54").Value = 19081
$ws.Range("J54").Value = 19081
$ws.Range("L54").Value = 19081
$ws.Range("N54").Value = -20369

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 3502
$ws.Range("I14").Value = 3502
$ws.Range("K14").Value = 3502
$ws.Range("M14").Value = -3334

$ws.Range("H28").Value = 80000
$ws.Range("J28").Value = 80000
$ws.Range("L28").Value = 80000
$ws.Range("N28").Value = -80696

$ws.Range("H29").Value = 35671
$ws.Range("I29").Value = 35671
$ws.Range("K29").Value = 35671
$ws.Range("M29").Value = -35381

$ws.Range("H51").Value = 25514.715
$ws.Range("I51").Value = 25228
$ws.Range("J51").Value = 25629.4
$ws.Range("K51").Value = 25228
$ws.Range("L51").Value = 25629.4
$ws.Range("M51").Value = -24718
$ws.Range("N51").Value = -26649.4

$ws.Range("H55").Value = 27166
$ws.Range("J55").Value = 33249.75
$ws.Range("L55").Value = 33249.75
$ws.Range("N55").Value = -33803.75

$ws.Range("H59").Value = 19999
$ws.Range("J59").Value = 19999
$ws.Range("L59").Value = 19999
$ws.Range("N59").Value = -21475
